$d = $word.ActiveDocument

# 1. Insert a new cue paragraph "Petra (neutral curious):" right before the
#    paragraph "I spin around and come face to face with Petra."
$r = $d.Content
$found = $r.Find.Execute("I spin around and come face to face with Petra.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.InsertParagraphBefore()
    $newIndex = $r.Paragraphs(1).Index
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = "Petra (neutral curious):"
}

# 2. Drop the stage-direction tag from the next line so it simply reads
#    "Petra: It's a little odd to see someone walking around like that."
$replaced2 = $d.Content.Find.Execute("Petra (neutral curious): It’s a little odd to see someone walking around like that.", $true, $false, $false, $false, $false, $true, 1, $false, "Petra: It’s a little odd to see someone walking around like that.", 2)

# 3. Remove the stray trailing empty run after "Petra (neutral confused):"
#    by rebuilding that paragraph cleanly.
$r3 = $d.Content
$found3 = $r3.Find.Execute("Petra (neutral confused):", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $pIndex = $r3.Paragraphs(1).Index
    $p = $d.Paragraphs.Item($pIndex)
    $p.Range.Delete()
    $after = $d.Paragraphs.Item($pIndex)
    $after.Range.InsertParagraphBefore()
    $rebuilt = $d.Paragraphs.Item($pIndex)
    $rebuilt.Range.Text = "Petra (neutral confused):"
}
